# Hightodown Doku.docx - add documentation for three new functions
# (random_in_alias, only_alias_else_jump, only_alias_endswith) right
# after the existing "only_alias" entry, and drop one of the two blank
# trailing paragraphs (the _GoBack bookmark moves to the end of the new
# content, as it does whenever Word's cursor is left at the end of the
# last edit).

$d = $word.ActiveDocument

# --- locate the paragraph whose text ends the "only_alias" entry -----
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Wert im Alias dem in der n*chsten Zeile entspricht*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "anchor paragraph not found"
}

# --- the _GoBack bookmark currently sits at the end of that paragraph;
#     remove it, it will be re-created at the end of the new content --
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- create an empty paragraph right after the anchor, to receive the
#     new OOXML content ---------------------------------------------
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$insertionPoint = $anchorPara.Range
$insertionPoint.Collapse(0)            # wdCollapseEnd
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$target = $newPara.Range

$body = '<w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:highlight w:val="darkYellow"/></w:rPr><w:t>random_in_alias</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Generiert eine Pseudozufallszahl die zwischen den in den n&#228;chsten beiden Zeilen angegebenen Werten und speichert sie dann im Alias</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:highlight w:val="darkYellow"/></w:rPr><w:t>only_alias_else_jump</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">code wird nur weiter </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>ausgef&#252;hrt</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> wenn der Alias dem entspricht was in der n&#228;chsten Zeile angegeben ist. Ansonsten wird in die </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Zeile gesprungen</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> die in der &#252;bern&#228;chsten angegeben wird</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/><w:highlight w:val="darkYellow"/></w:rPr><w:t>only_alias_endswith</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Code wird nur weiter ausgef&#252;hrt wenn der Alias mit dem in der n&#228;chsten Zeile angegebenen Zeichen endet</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)

# --- the document used to end in two empty paragraphs; now that the
#     bookmark has moved onto real content, drop the first of the two
#     so only one blank trailing paragraph remains, as in the target --
# Walk backwards from the end of the document while paragraphs are
# empty (just a pilcrow) and remember the topmost one of that run.
$count2 = $d.Paragraphs.Count
$i = $count2
$emptyRunStart = -1
while ($i -ge 1) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq [char]13) {
        $emptyRunStart = $i
        $i = $i - 1
    } else {
        break
    }
}
if ($emptyRunStart -ne -1 -and $emptyRunStart -lt $count2) {
    $d.Paragraphs.Item($emptyRunStart).Range.Delete()
}

Write-Output "done"
